$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: date, tag, amount (now text), description
$ws.Range("B3").Value = 241101
$ws.Range("C3").Value = "bills"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "12"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "edited"

# Delete row 6 entirely (shift cells up)
$ws.Rows(6).Delete()
